$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference")

# "Number of Globals" (I2, named _1st_ref) dropped from 16 to 14 after removing
# the molten-salt parameter block elsewhere in the model. K2:K7 (which chain off
# I2/_1st_ref and Param_Count) recalc automatically.
$ws.Range("I2").Value = 14

# Row heights shifted: row 2 (the row that used to carry the custom 14.65pt
# height) drops back to the sheet's default height, while rows 5-15 pick up
# an explicit 14.65pt height.
$ws.Rows.Item(2).AutoFit()
$ws.Range("A5:A15").EntireRow.RowHeight = 14.65

# Active cell/selection moved from N5 to I3.
$ws.Range("I3").Select()
